# Conditional distributions section now complete.
#
# Adds the "p(X_B|X_A=1)" / "p(X_A)" labels, highlights the two posterior
# inputs (D3, D4) in red, and swaps the SUM() helper formulas in E3/E4 for
# the worked-out Bayes-rule arithmetic as descriptive text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New column header: p(X_B|X_A=1) --------------------------------------
$ws.Range("E2").Value = "p(X_B|X_A=1)"
[void]$ws.Range("B3").Copy()
[void]$ws.Range("E2").PasteSpecial(-4122)   # xlPasteFormats - match the bold/centered header look

# --- New row label: p(X_A) -------------------------------------------------
$ws.Range("B5").Value = "p(X_A)"
[void]$ws.Range("B3").Copy()
[void]$ws.Range("B5").PasteSpecial(-4122)   # xlPasteFormats

# --- Highlight the two posterior inputs in red ------------------------------
$ws.Range("D3:D4").Font.Color = 255   # RGB(255,0,0)

# --- Replace the SUM() formulas with the worked Bayes-rule text -----------
$ws.Range("E3").Value = '"=0.1/0.6 = 1/6"'
$ws.Range("E4").Value = '"=0.5/0.6 = 5/6"'

$excel.CutCopyMode = $false

# --- Cosmetics: widen the new column to fit its text, restore selection ---
[void]$ws.Columns("E").AutoFit()
[void]$ws.Range("I13").Select()
